# Updates cryptos list values (price + 1h volume change) to match the
# latest scrape, and restores the MXToken/Aave row order swap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value() = '29.483.33'
$ws.Range('E2').Value() = '  +0.75%  '
$ws.Range('D3').Value() = '1.913.35'
$ws.Range('E3').Value() = '  +0.11%  '
$ws.Range('E4').Value() = '  +0.60%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value() = '325.81'
$c.Style = "Normal"
$ws.Range('E5').Value() = '  +1.42%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value() = '1.007'
$c.Style = "Normal"
$ws.Range('E6').Value() = '  +0.52%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value() = '0.4819'
$c.Style = "Normal"
$ws.Range('E7').Value() = '  +2.03%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value() = '0.4066'
$c.Style = "Normal"
$ws.Range('E8').Value() = '  +0.04%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value() = '0.08152'
$c.Style = "Normal"
$ws.Range('E9').Value() = '  +1.41%  '
$ws.Range('E10').Value() = '  +1.10%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value() = '23.42'
$c.Style = "Normal"
$ws.Range('E11').Value() = '  +4.19%  '
$ws.Range('D12').Value() = '1.891.81'
$ws.Range('E12').Value() = '  -0.65%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value() = '6.014'
$c.Style = "Normal"
$ws.Range('E13').Value() = '  +2.06%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value() = '7.168'
$c.Style = "Normal"
$ws.Range('E14').Value() = '  +0.59%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value() = '90.24'
$c.Style = "Normal"
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value() = '0.06793'
$c.Style = "Normal"
$ws.Range('E16').Value() = '  +2.52%  '
$ws.Range('E17').Value() = '  +0.63%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value() = '17.71'
$c.Style = "Normal"
$ws.Range('E19').Value() = '  +0.15%  '
$ws.Range('E20').Value() = '  +0.49%  '
$ws.Range('D21').Value() = '29.504.63'
$ws.Range('E21').Value() = '  +0.75%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value() = '5.627'
$c.Style = "Normal"
$ws.Range('E22').Value() = '  +2.01%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value() = '2.183'
$c.Style = "Normal"
$ws.Range('E24').Value() = '  -0.62%  '
$ws.Range('D25').Value() = '2.144.51'
$ws.Range('E25').Value() = '  +0.63%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value() = '155.79'
$c.Style = "Normal"
$ws.Range('E26').Value() = '  +0.41%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value() = '6.374'
$c.Style = "Normal"
$ws.Range('E27').Value() = '  +6.22%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value() = '20.06'
$c.Style = "Normal"
$ws.Range('E28').Value() = '  +1.57%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value() = '2.107'
$c.Style = "Normal"
$ws.Range('E29').Value() = '  +0.20%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value() = '119.96'
$c.Style = "Normal"
$ws.Range('E30').Value() = '  +2.20%  '
$ws.Range('E31').Value() = '  -4.76%  '
$ws.Range('E32').Value() = '  +0.18%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value() = '5.519'
$c.Style = "Normal"
$ws.Range('E33').Value() = '  +2.50%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value() = '3.560'
$c.Style = "Normal"
$ws.Range('E34').Value() = '  +0.60%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value() = '1.391'
$c.Style = "Normal"
$ws.Range('E35').Value() = '  -2.10%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value() = '0.02269'
$c.Style = "Normal"
$ws.Range('E36').Value() = '  +1.14%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value() = '0.06104'
$c.Style = "Normal"
$ws.Range('E37').Value() = '  +0.45%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value() = '1.177'
$c.Style = "Normal"
$ws.Range('E38').Value() = '  +0.19%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value() = '0.5961'
$c.Style = "Normal"
$ws.Range('E39').Value() = '  +1.86%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value() = '10.80'
$c.Style = "Normal"
$ws.Range('E40').Value() = '  +6.91%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value() = '7.998'
$c.Style = "Normal"
$ws.Range('E41').Value() = '  -2.91%  '
$ws.Range('E42').Value() = '  +1.01%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value() = '1.280'
$c.Style = "Normal"
$ws.Range('E43').Value() = '  +0.73%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value() = '2.392'
$c.Style = "Normal"
$ws.Range('E44').Value() = '  -4.48%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value() = '12.50'
$c.Style = "Normal"
$ws.Range('E45').Value() = '  +3.09%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value() = '0.07605'
$c.Style = "Normal"
$ws.Range('E46').Value() = '  -3.63%  '
$ws.Range('E47').Value() = '  +0.87%  '
$ws.Range('E48').Value() = '  +0.90%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value() = '115.93'
$c.Style = "Normal"
$ws.Range('E49').Value() = '  +2.59%  '
$ws.Range('B50').Value() = 'Aave'
$ws.Range('C50').Value() = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value() = '72.53'
$c.Style = "Normal"
$ws.Range('E50').Value() = '  +1.82%  '
$ws.Range('B51').Value() = 'MXToken'
$ws.Range('C51').Value() = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value() = '2.415'
$c.Style = "Normal"
$ws.Range('E51').Value() = '  +2.74%  '
